$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 108: the original entry was logged as PM instead of AM; fix the
#     start/stop times to the correct (early-morning) values. ---
$ws.Range("B108").Value = 0.097222222222222224
$ws.Range("C108").Value = 0.15347222222222223

# --- Row 109: add a new time-log entry (was a blank filler row).
#     Set the "Interruption" minutes first so the shared formula in E109
#     picks up every precedent correctly when B/C are written. ---
$ws.Range("D109").Value = 10
$ws.Range("A109").Value = 41951
$ws.Range("B109").Value = 0.60416666666666663
$ws.Range("C109").Value = 0.67361111111111116
$ws.Range("F109").Value = "Coding"

# Recalculate the workbook so Sheet2's SUMIF/percentage formulas and the
# pie chart cache pick up the new totals.
$excel.CalculateFull()

# --- Restore the window scroll position / selection on Sheet1 ---
$excel.ActiveWindow.ScrollRow = 97
$ws.Range("C110").Select()
